# papers_duedates.xlsx - "Revised 2019 charts. New 2020 charts."
#
# The workbook has 3 sheets; the edited data lives on "due_dates_2019jun"
# (index 2). The edit:
#   1. Revises several 2019 due_date values (column C) on rows 2-6.
#   2. Adds a new column K (rows 5-16) holding a month-stepped running
#      series of serial numbers (the "2020 chart" helper series), built
#      with formulas that add on each month's day-count to the prior cell.
#   3. Updates the sheet's view (zoom + selected cell) to reflect where the
#      user ended up working (K16), since the used range grew to column K.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("due_dates_2019jun")
$ws.Activate()

# --- 1. Revised 2019 due dates (column C, rows 2-6) -------------------
$ws.Range("C2").Value2 = 43466
$ws.Range("C3").Value2 = 43480
$ws.Range("C4").Value2 = 43496
$ws.Range("C5").Value2 = 43506
$ws.Range("C6").Value2 = 43521

# --- 2. New column K: month-stepped series for rows 5-16 --------------
$ws.Range("K5").Value2 = 21915

$ws.Range("K6").Formula = "=K5+31"
$ws.Range("K7").Formula = "=K6+28"
$ws.Range("K8").Formula = "=K7+31"
$ws.Range("K9").Formula = "=K8+30"
$ws.Range("K10").Formula = "=K9+31"
$ws.Range("K11").Formula = "=K10+30"
$ws.Range("K12").Formula = "=K11+31"
$ws.Range("K13").Formula = "=K12+31"
$ws.Range("K14").Formula = "=K13+30"
$ws.Range("K15").Formula = "=K14+31"
$ws.Range("K16").Formula = "=K15+31"

# --- 3. View state: zoom + active selection moved to the new cell -----
[void]$ws.Range("K16").Select()
$excel.ActiveWindow.Zoom = 115
